$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 8102.3687
$ws.Range("I28").Value = 233.6875
$ws.Range("J28").Value = 50068.668
$ws.Range("K28").Value = 233.6875
$ws.Range("L28").Value = 50068.668
$ws.Range("M28").Value = 251.3125
$ws.Range("N28").Value = -51038.668

# Row 40
$ws.Range("H40").Value = 3685.8
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 4107.25
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 4107.25
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -4457.25

# Row 62
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5248

# Row 65
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -26240

# Row 98
$ws.Range("H98").Value = 1035.3334
$ws.Range("I98").Value = 503
$ws.Range("K98").Value = 503
$ws.Range("M98").Value = 995

# Row 113
$ws.Range("H113").Value = 2425
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 2850
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 2850
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -9358

# Row 122
$ws.Range("H122").Value = 1035.3334
$ws.Range("I122").Value = 503
$ws.Range("K122").Value = 1509
$ws.Range("M122").Value = 941

# Row 132
$ws.Range("H132").Value = 4361.28
$ws.Range("I132").Value = 3896.4736
$ws.Range("J132").Value = 5833.1665
$ws.Range("K132").Value = 11689.4208
$ws.Range("L132").Value = 17499.4995
$ws.Range("M132").Value = -9159.4208
$ws.Range("N132").Value = -22559.4995

# Row 138
$ws.Range("H138").Value = 1830.1
$ws.Range("I138").Value = 572.8421
$ws.Range("J138").Value = 2125.0125
$ws.Range("K138").Value = 1718.5263
$ws.Range("L138").Value = 6375.037499999999
$ws.Range("M138").Value = 3421.4737
$ws.Range("N138").Value = -16655.0375

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 394968.5
$ws.Range("I32").Value = 465369.47
$ws.Range("J32").Value = 15886.462
$ws.Range("K32").Value = 465369.47
$ws.Range("L32").Value = 15886.462
$ws.Range("M32").Value = -465082.47
$ws.Range("N32").Value = -16460.462

# Row 37
$ws.Range("H37").Value = 13999.333

# Row 44
$ws.Range("H44").Value = 18623.5
$ws.Range("J44").Value = 18623.5
$ws.Range("L44").Value = 18623.5
$ws.Range("N44").Value = -19599.5

# Row 45
$ws.Range("H45").Value = 2710.3635
$ws.Range("I45").Value = 2441.4
$ws.Range("J45").Value = 2934.5
$ws.Range("K45").Value = 2441.4
$ws.Range("L45").Value = 2934.5
$ws.Range("M45").Value = -2064.4
$ws.Range("N45").Value = -3688.5

# Row 63
$ws.Range("H63").Value = 5137.7856
$ws.Range("I63").Value = 2992.875
$ws.Range("K63").Value = 2992.875
$ws.Range("M63").Value = -2306.875

# Row 66
$ws.Range("H66").Value = 5137.7856
$ws.Range("I66").Value = 2992.875
$ws.Range("K66").Value = 14964.375
$ws.Range("M66").Value = -11532.375

# Row 122
$ws.Range("H122").Value = 85300
$ws.Range("I122").Value = 112600
$ws.Range("J122").Value = 3400
$ws.Range("K122").Value = 337800
$ws.Range("L122").Value = 10200
$ws.Range("M122").Value = -335350
$ws.Range("N122").Value = -15100

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1490.9354
$ws.Range("I22").Value = 1551.1666
$ws.Range("J22").Value = 1284.4286
$ws.Range("K22").Value = 1551.1666
$ws.Range("L22").Value = 1284.4286
$ws.Range("M22").Value = -1378.1666
$ws.Range("N22").Value = -1630.4286

# Row 82
$ws.Range("H82").Value = 14771.857
$ws.Range("J82").Value = 21249
$ws.Range("L82").Value = 21249
$ws.Range("N82").Value = -22015

# Row 85
$ws.Range("H85").Value = 14771.857
$ws.Range("J85").Value = 21249
$ws.Range("L85").Value = 21249
$ws.Range("N85").Value = -23901

# Row 99
$ws.Range("H99").Value = 1879.8
$ws.Range("I99").Value = 1933
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1933
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -435
$ws.Range("N99").Value = -4796

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1932.4546
$ws.Range("I99").Value = 1900
$ws.Range("J99").Value = 1934
$ws.Range("K99").Value = 1900
$ws.Range("L99").Value = 1934
$ws.Range("M99").Value = -402
$ws.Range("N99").Value = -4930

# Row 108
$ws.Range("H108").Value = 26332.75
$ws.Range("J108").Value = 26332.75
$ws.Range("L108").Value = 26332.75
$ws.Range("N108").Value = -34012.75

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 112
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954

# Row 126
$ws.Range("H126").Value = 1932.4546
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 1934
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 5802
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -10742

$ws = $wb.Worksheets.Item("CUL")
# Row 119
$ws.Range("H119").Value = 3101.8
$ws.Range("I119").Value = 2627.25
$ws.Range("K119").Value = 7881.75
$ws.Range("M119").Value = -3043.75

# Row 122
$ws.Range("H122").Value = 9671.637000000001
$ws.Range("I122").Value = 348.75
$ws.Range("J122").Value = 34532.668
$ws.Range("K122").Value = 3138.75
$ws.Range("L122").Value = 310794.012
$ws.Range("M122").Value = -688.75
$ws.Range("N122").Value = -315694.012

# Row 125
$ws.Range("H125").Value = 2203.3333
$ws.Range("J125").Value = 2203.3333
$ws.Range("L125").Value = 6609.999899999999
$ws.Range("N125").Value = -16449.9999

# Row 131
$ws.Range("H131").Value = 970.4872
$ws.Range("J131").Value = 1041.9714
$ws.Range("L131").Value = 3125.9142
$ws.Range("N131").Value = -13205.9142

# Row 139
$ws.Range("H139").Value = 2623.2144
$ws.Range("I139").Value = 1931.875
$ws.Range("J139").Value = 2899.75
$ws.Range("K139").Value = 5795.625
$ws.Range("L139").Value = 8699.25
$ws.Range("M139").Value = -655.625
$ws.Range("N139").Value = -18979.25

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 251153.25
$ws.Range("J113").Value = 1537.6666
$ws.Range("L113").Value = 1537.6666
$ws.Range("N113").Value = -5877.6666

# Row 122
$ws.Range("H122").Value = 3862.5833
$ws.Range("I122").Value = 2801.4
$ws.Range("J122").Value = 4620.5713
$ws.Range("K122").Value = 8404.200000000001
$ws.Range("L122").Value = 13861.7139
$ws.Range("M122").Value = -5954.200000000001
$ws.Range("N122").Value = -18761.7139

# Row 132
$ws.Range("H132").Value = 3252.375
$ws.Range("I132").Value = 2776.5
$ws.Range("J132").Value = 3728.25
$ws.Range("K132").Value = 8329.5
$ws.Range("L132").Value = 11184.75
$ws.Range("M132").Value = -5799.5
$ws.Range("N132").Value = -16244.75

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 61203.6
$ws.Range("J20").Value = 61203.6
$ws.Range("L20").Value = 61203.6
$ws.Range("N20").Value = -61655.6

# Row 40
$ws.Range("H40").Value = 1000004
$ws.Range("I40").Value = 1000004
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1000004
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -999868
$ws.Range("N40").ClearContents()

# Row 46
$ws.Range("H46").Value = 4308.3335
$ws.Range("I46").Value = 1050
$ws.Range("J46").Value = 5937.5
$ws.Range("K46").Value = 1050
$ws.Range("L46").Value = 5937.5
$ws.Range("M46").Value = -862
$ws.Range("N46").Value = -6313.5

# Row 122
$ws.Range("H122").Value = 2965.5334
$ws.Range("I122").Value = 2857.5483
$ws.Range("J122").Value = 3204.6428
$ws.Range("K122").Value = 8572.644899999999
$ws.Range("L122").Value = 9613.928400000001
$ws.Range("M122").Value = -6122.644899999999
$ws.Range("N122").Value = -14513.9284

# Row 136
$ws.Range("H136").Value = 9262008
$ws.Range("I136").Value = 3750
$ws.Range("J136").Value = 11907224
$ws.Range("K136").Value = 11250
$ws.Range("L136").Value = 35721672
$ws.Range("M136").Value = -8700
$ws.Range("N136").Value = -35726772

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 87535.53
$ws.Range("I29").Value = 9363.637000000001
$ws.Range("J29").Value = 302508.25
$ws.Range("K29").Value = 9363.637000000001
$ws.Range("L29").Value = 302508.25
$ws.Range("M29").Value = -9073.637000000001
$ws.Range("N29").Value = -303088.25

# Row 122
$ws.Range("H122").Value = 4091.4285
$ws.Range("I122").Value = 3320
$ws.Range("J122").Value = 4670
$ws.Range("K122").Value = 9960
$ws.Range("L122").Value = 14010
$ws.Range("M122").Value = -7510
$ws.Range("N122").Value = -18910

# Row 136
$ws.Range("H136").Value = 4627.7856
$ws.Range("I136").Value = 4498.9
$ws.Range("J136").Value = 4950
$ws.Range("K136").Value = 13496.7
$ws.Range("L136").Value = 14850
$ws.Range("M136").Value = -10946.7
$ws.Range("N136").Value = -19950
